# Add the new worksheet "analissi inferenziale" as the last sheet in the workbook.
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "analissi inferenziale"

# --- Header row (bold) ---
$ws.Range("B4").Value = "Statistica"
$ws.Range("C4").Value = "Gruppo 1 (n=8)"
$ws.Range("D4").Value = "Gruppo 2 (n=6)"

# --- Corrette (mean +/- SD) ---
$ws.Range("B5").Value = "Corrette (mean ± SD)"
$ws.Range("C5").Value = "13.5 ± 3.07"
$ws.Range("D5").Value = "18.7 ± 1.03"

# --- Sbagliate (mean +/- SD) ---
$ws.Range("B6").Value = "Sbagliate (mean ± SD)"
$ws.Range("C6").Value = "6.5 ± 3.07"
$ws.Range("D6").Value = "1.33 ± 1.03"

# --- Totale domande ---
$ws.Range("B7").Value = "Totale domande"
$ws.Range("C7").Value = 20
$ws.Range("D7").Value = 20

# --- % corrette (mean +/- SD) ---
$ws.Range("B8").Value = "% corrette (mean ± SD)"
$ws.Range("C8").Value = "67.5 ± 15.35"
$ws.Range("D8").Value = "93.3 ± 5.16"

# --- Min % corrette ---
$ws.Range("B9").Value = "Min % corrette"
$ws.Range("C9").Value = 45
$ws.Range("D9").Value = 85

# --- Max % corrette ---
$ws.Range("B10").Value = "Max % corrette"
$ws.Range("C10").Value = 90
$ws.Range("D10").Value = 100

# --- 25 percentile % corrette (values stored as text "63.75"/"91.25") ---
$ws.Range("B11").Value = "25° percentile % corrette"
$ws.Range("C11").Formula = '="63.75"'
$ws.Range("D11").Formula = '="91.25"'
$ws.Range("C11:D11").Copy() | Out-Null
$ws.Range("C11:D11").PasteSpecial(-4163) | Out-Null

# --- 50 percentile % corrette ---
$ws.Range("B12").Value = "50° percentile % corrette"
$ws.Range("C12").Value = 70
$ws.Range("D12").Value = 95

# --- 75 percentile % corrette ---
$ws.Range("B13").Value = "75° percentile % corrette"
$ws.Range("C13").Value = 75
$ws.Range("D13").Value = 95

# --- Section title ---
$ws.Range("B15").Value = "Test statistici tra i gruppi:"

# --- Test results table header (bold) ---
$ws.Range("B17").Value = "Test"
$ws.Range("C17").Value = "Valore"
$ws.Range("D17").Value = "p-value"
$ws.Range("E17").Value = "Interpretazione"

# --- T-test row ---
$ws.Range("B18").Value = "T-test (media % corrette)"
$ws.Range("C18").Value = "t = -4.436"
$ws.Range("D18").Formula = '="0.0016"'
$ws.Range("D18").Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4163) | Out-Null
$ws.Range("E18").Value = "Differenza significativa tra i gruppi"

# --- Chi-square row ---
$ws.Range("B19").Value = "Chi-quadro corrette/errate"
$ws.Range("C19").Value = "χ² = 25.667"
$ws.Range("D19").Value = "< 0.0001"
$ws.Range("E19").Value = "Distribuzione risposte significativamente diversa"

# --- Bold styling for header / label cells ---
$ws.Range("B4:D4").Font.Bold = $true
$ws.Range("B5").Font.Bold = $true
$ws.Range("B6").Font.Bold = $true
$ws.Range("B7").Font.Bold = $true
$ws.Range("B8").Font.Bold = $true
$ws.Range("B9").Font.Bold = $true
$ws.Range("B10").Font.Bold = $true
$ws.Range("B11").Font.Bold = $true
$ws.Range("B12").Font.Bold = $true
$ws.Range("B13").Font.Bold = $true
$ws.Range("B15").Font.Bold = $true
$ws.Range("B17:E17").Font.Bold = $true

# --- Selection matching the saved state ---
$ws.Range("D4").Select() | Out-Null

Write-Output "analissi inferenziale sheet added"
